$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 23.19000000000019
$ws.Range("H2").Value = 0.01343599788550887
$ws.Range("I2").Value = 0.01343599788550887
$ws.Range("L2").Value = 41.87606698797467
$ws.Range("M2").Value = '[6.955798295353759, 76.79633568059559]'
$ws.Range("N2").Value = 0.01984665625320914
$ws.Range("O2").Value = 0.01984665625320914
$ws.Range("P2").Value = 1.17613178422681
$ws.Range("Q2").Value = '[0.044026323473730145, 2.308237244979889]'
$ws.Range("R2").Value = 0.04207097597411513
$ws.Range("S2").Value = 0.04207097597411513
$ws.Range("T2").Value = 68.00994989157491
$ws.Range("U2").Value = '[48.92375302260645, 87.09614676054336]'
$ws.Range("V2").Value = [double]"5.574315009582165e-09"
$ws.Range("W2").Value = [double]"5.574315009582165e-09"
$ws.Range("X2").Value = 18.84912912912928
$ws.Range("Y2").Value = 14.67075075075087
$ws.Range("Z2").Value = 23.0275075075077
$ws.Range("F3").Value = 23.19000000000019
$ws.Range("H3").Value = 0.001207879630190756
$ws.Range("I3").Value = 0.001207879630190756
$ws.Range("L3").Value = 48.78507077288676
$ws.Range("M3").Value = '[18.40212226116242, 79.1680192846111]'
$ws.Range("N3").Value = 0.002289239445157909
$ws.Range("O3").Value = 0.002289239445157909
$ws.Range("P3").Value = 1.742184514603348
$ws.Range("Q3").Value = '[0.9874475407679633, 2.4969214884387334]'
$ws.Range("R3").Value = [double]"2.928659974665848e-05"
$ws.Range("S3").Value = [double]"2.928659974665848e-05"
$ws.Range("T3").Value = 57.73644792889424
$ws.Range("U3").Value = '[40.43481596373044, 75.03807989405803]'
$ws.Range("V3").Value = [double]"2.643492580389761e-08"
$ws.Range("W3").Value = [double]"2.643492580389761e-08"
$ws.Range("X3").Value = 16.75993993994008
$ws.Range("Y3").Value = 13.97435435435447
$ws.Range("Z3").Value = 19.54552552552568
$ws.Range("F4").Value = 23.19000000000019
$ws.Range("H4").Value = 0.0001668106417727078
$ws.Range("I4").Value = 0.0001668106417727078
$ws.Range("L4").Value = 47.881077565312
$ws.Range("M4").Value = '[19.10531112395506, 76.65684400666893]'
$ws.Range("N4").Value = 0.001636756385126725
$ws.Range("O4").Value = 0.001636756385126725
$ws.Range("P4").Value = 2.207605648468504
$ws.Range("Q4").Value = '[1.628973968528042, 2.7862373284089657]'
$ws.Range("R4").Value = [double]"9.970331227293627e-10"
$ws.Range("S4").Value = [double]"9.970331227293627e-10"
$ws.Range("T4").Value = 62.16368245085517
$ws.Range("U4").Value = '[47.41766662861126, 76.90969827309908]'
$ws.Range("V4").Value = [double]"6.730660473408534e-11"
$ws.Range("W4").Value = [double]"6.730660473408534e-11"
$ws.Range("X4").Value = 15.04216216216228
$ws.Range("Y4").Value = 12.90654654654665
$ws.Range("Z4").Value = 17.17777777777791
$ws.Range("F5").Value = 23.19000000000019
$ws.Range("H5").Value = [double]"4.030664371201809e-05"
$ws.Range("I5").Value = [double]"4.030664371201809e-05"
$ws.Range("L5").Value = 52.05164079903686
$ws.Range("M5").Value = '[23.60524261705595, 80.49803898101777]'
$ws.Range("N5").Value = 0.0006112132057289887
$ws.Range("O5").Value = 0.0006112132057289887
$ws.Range("P5").Value = 2.220184598032427
$ws.Range("Q5").Value = '[1.6918687163476571, 2.7485004797171966]'
$ws.Range("R5").Value = [double]"7.351341757555474e-11"
$ws.Range("S5").Value = [double]"7.351341757555474e-11"
$ws.Range("T5").Value = 65.69441373685737
$ws.Range("U5").Value = '[51.07335191270049, 80.31547556101425]'
$ws.Range("V5").Value = [double]"1.078381828278907e-11"
$ws.Range("W5").Value = [double]"1.078381828278907e-11"
$ws.Range("X5").Value = 14.99573573573585
$ws.Range("Y5").Value = 13.04582582582593
$ws.Range("Z5").Value = 16.94564564564578
$ws.Range("F6").Value = 23.19000000000019
$ws.Range("H6").Value = 0.0001473023574158905
$ws.Range("I6").Value = 0.0001473023574158905
$ws.Range("L6").Value = 52.41525418874421
$ws.Range("M6").Value = '[21.542498765057594, 83.28800961243081]'
$ws.Range("N6").Value = 0.001343420394590789
$ws.Range("O6").Value = 0.001343420394590789
$ws.Range("P6").Value = 2.320816194543811
$ws.Range("Q6").Value = '[1.7673424137311944, 2.8742899753564277]'
$ws.Range("R6").Value = [double]"7.816147729045042e-11"
$ws.Range("S6").Value = [double]"7.816147729045042e-11"
$ws.Range("T6").Value = 62.42318838088539
$ws.Range("U6").Value = '[46.368839339901164, 78.47753742186961]'
$ws.Range("V6").Value = [double]"6.07341288372254e-10"
$ws.Range("W6").Value = [double]"6.07341288372254e-10"
$ws.Range("X6").Value = 14.62432432432444
$ws.Range("Y6").Value = 12.58156156156166
$ws.Range("Z6").Value = 16.66708708708722
$ws.Range("F7").Value = 23.19000000000019
$ws.Range("H7").Value = 0.0172851074959286
$ws.Range("I7").Value = 0.0172851074959286
$ws.Range("L7").Value = 34.90142076106215
$ws.Range("M7").Value = '[2.5105072640055823, 67.29233425811871]'
$ws.Range("N7").Value = 0.0353080976400364
$ws.Range("O7").Value = 0.0353080976400364
$ws.Range("P7").Value = 2.295658295415965
$ws.Range("Q7").Value = '[1.490605523324886, 3.100711067507044]'
$ws.Range("R7").Value = [double]"7.496224598835255e-07"
$ws.Range("S7").Value = [double]"7.496224598835255e-07"
$ws.Range("T7").Value = 58.84084870443866
$ws.Range("U7").Value = '[42.183675545632624, 75.49802186324469]'
$ws.Range("V7").Value = [double]"6.888697168250246e-09"
$ws.Range("W7").Value = [double]"6.888697168250246e-09"
$ws.Range("X7").Value = 14.71717717717729
$ws.Range("Y7").Value = 11.74588588588598
$ws.Range("Z7").Value = 17.68846846846861
$ws.Range("F8").Value = 23.19000000000019
$ws.Range("H8").Value = 0.03002130504994938
$ws.Range("I8").Value = 0.03002130504994938
$ws.Range("L8").Value = 25.37034603148295
$ws.Range("M8").Value = '[1.395945669992237, 49.34474639297367]'
$ws.Range("N8").Value = 0.03855445760633258
$ws.Range("O8").Value = 0.03855445760633258
$ws.Range("P8").Value = 2.836553126664658
$ws.Range("Q8").Value = '[1.4906055233248865, 4.1825007300044295]'
$ws.Range("R8").Value = 0.0001080322191406857
$ws.Range("S8").Value = 0.0001080322191406857
$ws.Range("T8").Value = 57.84670108400231
$ws.Range("U8").Value = '[44.55061353371607, 71.14278863428854]'
$ws.Range("V8").Value = [double]"2.749711569549618e-11"
$ws.Range("W8").Value = [double]"2.749711569549618e-11"
$ws.Range("X8").Value = 12.72084084084094
$ws.Range("Y8").Value = 7.753213213213273
$ws.Range("Z8").Value = 17.68846846846861
$ws.Range("F9").Value = 22.80000000000013
$ws.Range("H9").Value = 0.0001671336584448957
$ws.Range("I9").Value = 0.0001671336584448957
$ws.Range("L9").Value = 43.70244293397591
$ws.Range("M9").Value = '[19.761269185095756, 67.64361668285606]'
$ws.Range("N9").Value = 0.0006277488214474491
$ws.Range("O9").Value = 0.0006277488214474491
$ws.Range("P9").Value = -2.930895248394081
$ws.Range("Q9").Value = '[-3.559842726590236, -2.3019477701979265]'
$ws.Range("R9").Value = [double]"3.645972412869014e-12"
$ws.Range("S9").Value = [double]"3.645972412869014e-12"
$ws.Range("T9").Value = 50.99527175085864
$ws.Range("U9").Value = '[37.22364215831597, 64.76690134340132]'
$ws.Range("V9").Value = [double]"2.143219823835807e-09"
$ws.Range("W9").Value = [double]"2.143219823835807e-09"
$ws.Range("X9").Value = 10.63543543543549
$ws.Range("Y9").Value = 8.353153153153201
$ws.Range("Z9").Value = 12.91771771771779
$ws.Range("F10").Value = 22.80000000000013
$ws.Range("H10").Value = 0.002634963855456363
$ws.Range("I10").Value = 0.002634963855456363
$ws.Range("L10").Value = 42.30841420109132
$ws.Range("M10").Value = '[12.372354217282478, 72.24447418490016]'
$ws.Range("N10").Value = 0.006634934540027393
$ws.Range("O10").Value = 0.006634934540027393
$ws.Range("P10").Value = -2.767368904063081
$ws.Range("Q10").Value = '[-3.5850006257180818, -1.9497371824080796]'
$ws.Range("R10").Value = [double]"1.904710766886808e-08"
$ws.Range("S10").Value = [double]"1.904710766886808e-08"
$ws.Range("T10").Value = 58.35811450640675
$ws.Range("U10").Value = '[41.777323573900965, 74.93890543891254]'
$ws.Range("V10").Value = [double]"7.523863532199471e-09"
$ws.Range("W10").Value = [double]"7.523863532199471e-09"
$ws.Range("X10").Value = 10.0420420420421
$ws.Range("Y10").Value = 7.075075075075114
$ws.Range("Z10").Value = 13.00900900900908
$ws.Range("B11").Value = 0
$ws.Range("F11").Value = 22.80000000000013
$ws.Range("H11").Value = 0.03751753348472819
$ws.Range("I11").Value = 0.03751753348472819
$ws.Range("L11").Value = 33.20466831919682
$ws.Range("M11").Value = '[-0.24151205687279287, 66.65084869526643]'
$ws.Range("N11").Value = 0.05160698475077963
$ws.Range("O11").Value = 0.05160698475077963
$ws.Range("P11").Value = -2.188737224122619
$ws.Range("Q11").Value = '[-4.06300070914716, -0.3144737390980774]'
$ws.Range("R11").Value = 0.02310638813359889
$ws.Range("S11").Value = 0.02310638813359889
$ws.Range("T11").Value = 68.99506441598379
$ws.Range("U11").Value = '[51.41464397004272, 86.57548486192486]'
$ws.Range("V11").Value = [double]"4.750020377031205e-10"
$ws.Range("W11").Value = [double]"4.750020377031205e-10"
$ws.Range("X11").Value = 7.942342342342387
$ws.Range("Y11").Value = 1.141141141141148
$ws.Range("Z11").Value = 14.74354354354363
$ws.Range("F12").Value = 22.80000000000013
$ws.Range("H12").Value = [double]"2.637604060717891e-05"
$ws.Range("I12").Value = [double]"2.637604060717891e-05"
$ws.Range("L12").Value = 55.22373023100331
$ws.Range("M12").Value = '[26.486083360859766, 83.96137710114685]'
$ws.Range("N12").Value = 0.0003482389193374402
$ws.Range("O12").Value = 0.0003482389193374402
$ws.Range("P12").Value = -2.088105627611234
$ws.Range("Q12").Value = '[-2.6164215092960044, -1.5597897459264631]'
$ws.Range("R12").Value = [double]"3.935094472495848e-10"
$ws.Range("S12").Value = [double]"3.935094472495848e-10"
$ws.Range("T12").Value = 60.02561227870063
$ws.Range("U12").Value = '[45.087311604601126, 74.96391295280014]'
$ws.Range("V12").Value = [double]"2.523861120096171e-10"
$ws.Range("W12").Value = [double]"2.523861120096171e-10"
$ws.Range("X12").Value = 7.577177177177218
$ws.Range("Y12").Value = 5.660060060060088
$ws.Range("Z12").Value = 9.494294294294347
